$d = $word.ActiveDocument

$replacements = @(
    @("2026-02-22 Sunday", "2026-02-23 Monday"),
    @("674÷7=", "511÷3="),
    @("218÷4=", "396÷4="),
    @("777÷9=", "746÷8="),
    @("777÷4=", "504÷7="),
    @("748÷6=", "307÷5="),
    @("847÷4=", "613÷6="),
    @("428÷7=", "132÷5="),
    @("646÷2=", "496÷7="),
    @("188÷2=", "549÷3="),
    @("154÷6=", "486÷8="),
    @("242÷4=", "452÷9="),
    @("314÷6=", "744÷4="),
    @("306÷4=", "707÷2="),
    @("855÷6=", "422÷5="),
    @("329÷7=", "482÷8="),
    @("583÷8=", "596÷9="),
    @("849÷6=", "617÷9="),
    @("741÷9=", "249÷4="),
    @("662÷4=", "199÷4="),
    @("196÷6=", "486÷2="),
    @("880÷9=", "105÷7="),
    @("860÷2=", "619÷3="),
    @("666÷7=", "329÷4="),
    @("626÷7=", "296÷3="),
    @("364÷4=", "214÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
